$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nightly "Updated cryptos list" refresh: new prices / 1h volume deltas for
# every ranked coin, plus a few rows where the ranking reshuffled (a coin
# moved up/down a spot, or EnergySwap entered the list in place of the coin
# that fell out of the top 50).
#
# Price/volume cells are stored as text (e.g. "30.016.40", "1.001",
# "  +0.16%  ") rather than numbers, so every write is prefixed with a
# leading apostrophe (Excel's "treat as text" quote-prefix) and the cell
# style is put back to Normal afterwards so no numeric formatting sticks.
function Set-TextValue($a1, $value) {
    $range = $ws.Range($a1)
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Row 2: Bitcoin
Set-TextValue "D2" "29.952.97"
Set-TextValue "E2" "  +0.08%  "

# Row 3: Ethereum
Set-TextValue "D3" "1.908.79"
Set-TextValue "E3" "  +0.74%  "

# Row 4: TetherUSD
Set-TextValue "D4" "1.001"
Set-TextValue "E4" "  -0.09%  "

# Row 5: XRP
Set-TextValue "D5" "0.7963"
Set-TextValue "E5" "  +6.00%  "

# Row 6: BNB
Set-TextValue "D6" "241.63"
Set-TextValue "E6" "  +0.80%  "

# Row 7: USDC
Set-TextValue "D7" "1.001"
Set-TextValue "E7" "  -0.06%  "

# Row 8: Cardano
Set-TextValue "D8" "0.3158"
Set-TextValue "E8" "  +3.70%  "

# Row 9: Solana
Set-TextValue "D9" "26.28"
Set-TextValue "E9" "  +5.48%  "

# Row 10: Dogecoin
Set-TextValue "D10" "0.06910"
Set-TextValue "E10" "  +1.29%  "

# Row 11: TRON
Set-TextValue "D11" "0.07996"
Set-TextValue "E11" "  +0.35%  "

# Row 12: WrappedEther
Set-TextValue "D12" "1.909.40"
Set-TextValue "E12" "  +0.68%  "

# Row 13: Polygon
Set-TextValue "D13" "0.7418"
Set-TextValue "E13" "  -0.61%  "

# Row 14: Polkadot
Set-TextValue "D14" "5.190"
Set-TextValue "E14" "  +0.04%  "

# Row 15: Litecoin
Set-TextValue "D15" "92.89"
Set-TextValue "E15" "  +1.99%  "

# Row 16: WrappedBTC
Set-TextValue "D16" "29.964.18"
Set-TextValue "E16" "  +0.08%  "

# Row 17: Avalanche
Set-TextValue "D17" "13.97"
Set-TextValue "E17" "  +0.93%  "

# Row 18: Uniswap
Set-TextValue "D18" "5.867"
Set-TextValue "E18" "  -3.58%  "

# Row 19: BitcoinCash
Set-TextValue "D19" "245.91"
Set-TextValue "E19" "  +4.92%  "

# Row 20: ShibaInu
Set-TextValue "D20" "0.000007744"
Set-TextValue "E20" "  +1.13%  "

# Row 21: Dai
Set-TextValue "D21" "1.001"
Set-TextValue "E21" "  -0.10%  "

# Row 22: WrappedliquidstakedEther2.0
Set-TextValue "D22" "2.154.99"
Set-TextValue "E22" "  +0.45%  "

# Row 23: BinanceUSD
Set-TextValue "D23" "1.001"
Set-TextValue "E23" "  -0.09%  "

# Row 24: Chainlink
Set-TextValue "D24" "6.827"
Set-TextValue "E24" "  -1.20%  "

# Row 25: Monero
Set-TextValue "D25" "167.61"
Set-TextValue "E25" "  +1.63%  "

# Row 26: Cosmos
Set-TextValue "D26" "9.222"
Set-TextValue "E26" "  +0.14%  "

# Row 27: Stellar
Set-TextValue "D27" "0.1391"
Set-TextValue "E27" "  +8.98%  "

# Row 28: EthereumClassic
Set-TextValue "D28" "18.93"
Set-TextValue "E28" "  +1.60%  "

# Row 29: LidoDAOToken
Set-TextValue "D29" "2.028"
Set-TextValue "E29" "  -0.49%  "

# Row 30: Toncoin
Set-TextValue "D30" "1.367"
Set-TextValue "E30" "  +2.33%  "

# Row 31: PancakeSwap
Set-TextValue "D31" "1.512"
Set-TextValue "E31" "  -0.01%  "

# Row 32: Filecoin
Set-TextValue "D32" "4.309"
Set-TextValue "E32" "  +1.02%  "

# Row 33: InternetComputer(DFINITY)
Set-TextValue "B33" "InternetComputer(DFINITY)"
Set-TextValue "C33" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D33" "4.087"
Set-TextValue "E33" "  +2.46%  "

# Row 34: Hedera
Set-TextValue "B34" "Hedera"
Set-TextValue "C34" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D34" "0.05522"
Set-TextValue "E34" "  +3.17%  "

# Row 35: ARBITRUM
Set-TextValue "D35" "1.255"
Set-TextValue "E35" "  +1.33%  "

# Row 36: ImmutableX
Set-TextValue "D36" "0.7300"
Set-TextValue "E36" "  +0.34%  "

# Row 37: HuobiToken
Set-TextValue "E37" "  +0.04%  "

# Row 38: VeChain
Set-TextValue "D38" "0.01921"
Set-TextValue "E38" "  -0.03%  "

# Row 39: MXToken
Set-TextValue "D39" "2.788"
Set-TextValue "E39" "  +0.90%  "

# Row 40: FraxShare
Set-TextValue "D40" "6.124"
Set-TextValue "E40" "  -1.14%  "

# Row 41: TheSandbox
Set-TextValue "D41" "0.4412"
Set-TextValue "E41" "  +0.19%  "

# Row 42: Aave
Set-TextValue "D42" "72.07"
Set-TextValue "E42" "  -0.42%  "

# Row 43: PaxDollar
Set-TextValue "D43" "1.001"
Set-TextValue "E43" "  -0.09%  "

# Row 44: TrustWalletToken
Set-TextValue "D44" "0.8312"
Set-TextValue "E44" "  +0.80%  "

# Row 45: RenderToken
Set-TextValue "D45" "1.872"
Set-TextValue "E45" "  -2.18%  "

# Row 46: Quant
Set-TextValue "D46" "100.76"
Set-TextValue "E46" "  -0.37%  "

# Row 47: Aptos
Set-TextValue "D47" "7.512"
Set-TextValue "E47" "  -0.47%  "

# Row 48: EnergySwap
Set-TextValue "B48" "EnergySwap"
Set-TextValue "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "9.724"
Set-TextValue "E48" "  -0.35%  "

# Row 49: Maker
Set-TextValue "B49" "Maker"
Set-TextValue "C49" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D49" "988.73"
Set-TextValue "E49" "  +7.96%  "

# Row 50: RocketPoolETH
Set-TextValue "B50" "RocketPoolETH"
Set-TextValue "C50" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D50" "2.063.39"
Set-TextValue "E50" "  +0.34%  "

# Row 51: Elrond
Set-TextValue "B51" "Elrond"
Set-TextValue "C51" "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue "D51" "36.28"
Set-TextValue "E51" "  +0.34%  "
